$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Tabla1" table holds ACI EPG static-binding rows. The schema/template
# naming scheme changed: SCH_BPD -> SCH_TEST, PROD-CO -> PROD-SITE-1,
# PROD-CO-CA -> PROD-STRETCH, and the vpc interface policy group was renamed
# from BPD_test_vpc_ipg -> Test_vpc_ipg. Update every data row accordingly;
# the calculated columns (full_epg / schema_site) recompute automatically.

# epg_schema_name column (D2:D6): SCH_BPD -> SCH_TEST
$ws.Range("D2").Value = "SCH_TEST"
$ws.Range("D3").Value = "SCH_TEST"
$ws.Range("D4").Value = "SCH_TEST"
$ws.Range("D5").Value = "SCH_TEST"
$ws.Range("D6").Value = "SCH_TEST"

# epg_template_name column (E2:E6): PROD-CO -> PROD-SITE-1 / PROD-CO-CA -> PROD-STRETCH
$ws.Range("E2").Value = "PROD-SITE-1"
$ws.Range("E3").Value = "PROD-SITE-1"
$ws.Range("E4").Value = "PROD-SITE-1"
$ws.Range("E5").Value = "PROD-STRETCH"
$ws.Range("E6").Value = "PROD-STRETCH"

# interface_policy_group column (M3, M6): BPD_test_vpc_ipg -> Test_vpc_ipg
$ws.Range("M3").Value = "Test_vpc_ipg"
$ws.Range("M6").Value = "Test_vpc_ipg"

# Re-fit columns P (full_epg) and Q (schema_site) now that the concatenated
# text is longer, matching the widened "best fit" columns from the saved file.
$ws.Columns.Item(16).ColumnWidth = 66
$ws.Columns.Item(17).ColumnWidth = 14.142857142857142

# Move the active selection, matching the author's last cursor position on save.
$ws.Range("P17").Select()
